$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 20.736841
$ws.Range("I11").Value = 20.736841
$ws.Range("K11").Value = 20.736841
$ws.Range("M11").Value = 119.263159

$ws.Range("H33").Value = 195.72728
$ws.Range("I33").Value = 123.42857
$ws.Range("K33").Value = 123.42857
$ws.Range("M33").Value = 105.57143

$ws.Range("H51").Value = 3698.7144
$ws.Range("I51").Value = 3697
$ws.Range("J51").Value = 3699
$ws.Range("K51").Value = 3697
$ws.Range("L51").Value = 3699
$ws.Range("M51").Value = -3213
$ws.Range("N51").Value = -4667

$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 3000
$ws.Range("N59").Value = -4114

$ws.Range("H62").Value = 6289
$ws.Range("J62").Value = 7998.3335
$ws.Range("L62").Value = 7998.3335
$ws.Range("N62").Value = -9246.3335

$ws.Range("H65").Value = 6289
$ws.Range("J65").Value = 7998.3335
$ws.Range("L65").Value = 39991.6675
$ws.Range("N65").Value = -46231.6675

$ws.Range("H101").Value = 25001680
$ws.Range("J101").Value = 2240.3333
$ws.Range("L101").Value = 6720.999899999999
$ws.Range("N101").Value = -9964.999899999999

$ws.Range("H105").Value = 33734.2
$ws.Range("J105").Value = 33734.2
$ws.Range("L105").Value = 33734.2
$ws.Range("N105").Value = -40722.2

$ws.Range("H137").Value = 1228.6428
$ws.Range("I137").Value = 1207.7693
$ws.Range("K137").Value = 3623.3079
$ws.Range("M137").Value = -1073.3079

$ws.Range("H138").Value = 2433.4211
$ws.Range("J138").Value = 2121.2
$ws.Range("L138").Value = 6363.599999999999
$ws.Range("N138").Value = -16643.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2086.4285
$ws.Range("I74").Value = 2376.5715
$ws.Range("K74").Value = 2376.5715
$ws.Range("M74").Value = -1502.5715

$ws.Range("H77").Value = 2086.4285
$ws.Range("I77").Value = 2376.5715
$ws.Range("K77").Value = 11882.8575
$ws.Range("M77").Value = -7514.8575

$ws.Range("H122").Value = 856.5
$ws.Range("I122").Value = 856.5
$ws.Range("K122").Value = 2569.5
$ws.Range("M122").Value = -119.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1970
$ws.Range("I31").Value = 1601.909
$ws.Range("J31").Value = 2307.4167
$ws.Range("K31").Value = 1601.909
$ws.Range("L31").Value = 2307.4167
$ws.Range("M31").Value = -1306.909
$ws.Range("N31").Value = -2897.4167

$ws.Range("H34").Value = 1970
$ws.Range("I34").Value = 1601.909
$ws.Range("J34").Value = 2307.4167
$ws.Range("K34").Value = 1601.909
$ws.Range("L34").Value = 2307.4167
$ws.Range("M34").Value = -1399.909
$ws.Range("N34").Value = -2711.4167

$ws.Range("H60").Value = 6263.875
$ws.Range("I60").Value = 5462.2
$ws.Range("J60").Value = 7600
$ws.Range("K60").Value = 5462.2
$ws.Range("L60").Value = 7600
$ws.Range("M60").Value = -4951.2
$ws.Range("N60").Value = -8622

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H109").Value = 59950
$ws.Range("J109").Value = 59950
$ws.Range("L109").Value = 59950
$ws.Range("N109").Value = -62030

$ws.Range("H122").Value = 1422
$ws.Range("I122").Value = 562.6667
$ws.Range("K122").Value = 1688.0001
$ws.Range("M122").Value = 761.9999

$ws.Range("H140").Value = 98888
$ws.Range("J140").Value = 98888
$ws.Range("L140").Value = 98888
$ws.Range("N140").Value = -109248

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7857623.5
$ws.Range("I4").Value = 7857623.5
$ws.Range("K4").Value = 23572870.5
$ws.Range("M4").Value = -23572758.5

$ws.Range("H60").Value = 634
$ws.Range("I60").Value = 477
$ws.Range("K60").Value = 1431
$ws.Range("M60").Value = -1180

$ws.Range("H122").Value = 400
$ws.Range("I122").Value = 100
$ws.Range("K122").Value = 900
$ws.Range("M122").Value = 1550

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9749.25
$ws.Range("I70").Value = 9398.6
$ws.Range("J70").Value = 9999.714
$ws.Range("K70").Value = 9398.6
$ws.Range("L70").Value = 9999.714
$ws.Range("M70").Value = -9128.6
$ws.Range("N70").Value = -10539.714

$ws.Range("H73").Value = 9749.25
$ws.Range("I73").Value = 9398.6
$ws.Range("J73").Value = 9999.714
$ws.Range("K73").Value = 9398.6
$ws.Range("L73").Value = 9999.714
$ws.Range("M73").Value = -8462.6
$ws.Range("N73").Value = -11871.714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 350
$ws.Range("J20").Value = 350
$ws.Range("L20").Value = 350
$ws.Range("N20").Value = -802

$ws.Range("H22").Value = 2479.2
$ws.Range("J22").Value = 2379.2856
$ws.Range("L22").Value = 2379.2856
$ws.Range("N22").Value = -2969.2856

$ws.Range("H27").Value = 2479.2
$ws.Range("J27").Value = 2379.2856
$ws.Range("L27").Value = 2379.2856
$ws.Range("N27").Value = -2593.2856

$ws.Range("H40").Value = 3005.1538
$ws.Range("I40").Value = 1709.7142
$ws.Range("K40").Value = 1709.7142
$ws.Range("M40").Value = -1573.7142

$ws.Range("H74").Value = 90000
$ws.Range("I74").Value = 90000
$ws.Range("K74").Value = 90000
$ws.Range("M74").Value = -89002

$ws.Range("H77").Value = 90000
$ws.Range("I77").Value = 90000
$ws.Range("K77").Value = 270000
$ws.Range("M77").Value = -265008

$ws.Range("H132").Value = 2376.577
$ws.Range("I132").Value = 2148.9546
$ws.Range("K132").Value = 6446.8638
$ws.Range("M132").Value = -3916.8638

$ws.Range("H139").Value = 110650
$ws.Range("I139").Value = 110650
$ws.Range("K139").Value = 110650
$ws.Range("M139").Value = -105510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1075000
$ws.Range("J2").Value = 150000
$ws.Range("L2").Value = 150000
$ws.Range("N2").Value = -150224

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H92").Value = 30050
$ws.Range("J92").Value = 30050
$ws.Range("L92").Value = 30050
$ws.Range("N92").Value = -35042

$ws.Range("H107").Value = 1131.5714
$ws.Range("J107").Value = 475
$ws.Range("L107").Value = 1425
$ws.Range("N107").Value = -5265

$ws.Range("H109").Value = 59000
$ws.Range("J109").Value = 59000
$ws.Range("L109").Value = 59000
$ws.Range("N109").Value = -61774

$ws.Range("H122").Value = 1299.4
$ws.Range("I122").Value = 1186.75
$ws.Range("K122").Value = 3560.25
$ws.Range("M122").Value = -1110.25

$ws.Range("H132").Value = 7586.75
$ws.Range("I132").Value = 8465.667
$ws.Range("J132").Value = 4950
$ws.Range("K132").Value = 25397.001
$ws.Range("L132").Value = 14850
$ws.Range("M132").Value = -22867.001
$ws.Range("N132").Value = -19910

$ws.Range("H136").Value = 4827.1
$ws.Range("I136").Value = 4713
$ws.Range("J136").Value = 4998.25
$ws.Range("K136").Value = 14139
$ws.Range("L136").Value = 14994.75
$ws.Range("M136").Value = -11589
$ws.Range("N136").Value = -20094.75
